$d = $word.ActiveDocument

# --- Section 1: "Reusability" bullet list (numId=1) ---
# Insert two new bullet items right after:
#   "Communication between the m5core2 and nucleo board ..." (paragraph 6)
$p = $d.Paragraphs(6)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs(7)
$np.Range.Text = "M5core2 and nucleo board communicate via UART, and the nucleo board and ultrasonic sensors do as well. UART communication is supported by almost all devices in this space, so reusing this will be easy " + [char]8211 + " provided that the nucleo board alternative supports 2 different UART channels."

$p = $d.Paragraphs(7)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs(8)
$np.Range.Text = "The turtlebot could be replaced by any autonomous robotic system, provided that they are able to run ROS " + [char]8211 + " which should be generic. This robot doesn" + [char]8217 + "t even necessarily need to have a LiDAR sensor, provided that the LiDAR sensor has (UART?) communication to the turtlebot."

# --- Section 2: "Reconfigurations or possible addition features" bullet list (numId=2) ---
# Insert two new bullet items right after:
#   "The nucleo board could be replaced by another board ..."
# (originally paragraph 10; now shifted by +2 because of the two paragraphs inserted above)
$p = $d.Paragraphs(12)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs(13)
$np.Range.Text = "Could combine multiple different sensors in order to broaden the gesture capabilities"

$p = $d.Paragraphs(13)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs(14)
$np.Range.Text = "Could replace sensor(s) with a camera and integrate machine learning in order to further diversify gestures."

# --- Section 3: "Ease of use" bullet list (numId=1) ---
# Insert two new bullet items right after:
#   "M5core2 only requires itself to be turned on ..."
# (originally paragraph 16; now shifted by +4 because of the four paragraphs inserted above)
$p = $d.Paragraphs(20)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs(21)
$np.Range.Text = "Manual describing gestures"

$p = $d.Paragraphs(21)
$p.Range.InsertParagraphAfter()
$np = $d.Paragraphs(22)
$np.Range.Text = "Simplify and outline wiring in a manual."

# One extra blank paragraph, inserted just before the document's pre-existing
# trailing blank paragraph (matches the diff's additional "<w:p/>").
$p = $d.Paragraphs(23)
$p.Range.InsertParagraphBefore()

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
